$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.165.02"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.883.67"
$ws.Range("E3").Value = "  +4.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.36%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").Value = "3.337.99"
$ws.Range("E15").Value = "  +4.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.996"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.05%  "
$ws.Range("D17").Value = "2.882.44"
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").Value = "52.167.30"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  +3.28%  "
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("E25").Value = "  +1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "53.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0940"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.10%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0460"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +7.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.17%  "
$ws.Range("D47").Value = "2.191.75"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("E48").Value = "  +5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +19.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.63%  "
